$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column F (dSF) values for several rows to match re-pulled / recalculated data
$ws.Range("F3").Value = -7
$ws.Range("F5").Value = -5
$ws.Range("F7").Value = -12
$ws.Range("F8").Value = 10
$ws.Range("F9").Value = 17
$ws.Range("F10").Value = -7
$ws.Range("F14").Value = 5
$ws.Range("F15").Value = 0
